$d = $word.ActiveDocument

# Locate the target paragraphs by their text rather than a hard-coded
# index, so the script is resilient to the exact paragraph numbering.
$queensPara = $null
$randomPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Queens - mated (y/n)*" -and $queensPara -eq $null) {
        $queensPara = $p
    }
    if ($t -like "*Randomness/variability to resource collection*" -and $randomPara -eq $null) {
        $randomPara = $p
    }
}

# 1. Remove the trailing "\" run after "Queens - mated (y/n)" (deleting
#    just the final run's single character so the other runs in the
#    paragraph stay untouched/unmerged).
$r7 = $queensPara.Range
$bs = $d.Range($r7.End - 2, $r7.End - 1)
$bs.Delete()

# 2. Apply strikethrough formatting to the "Randomness/variability to
#    resource collection" paragraph (sets both the paragraph mark's run
#    properties and the run's own properties, matching how Word toggles
#    Font.StrikeThrough over a paragraph-spanning range).
$randomPara.Range.Font.StrikeThrough = 1
